# Applies the "updated cryptos list" data refresh described by the diff.
# Values are written with a leading apostrophe (forced-text marker) so that
# numeric-looking strings (e.g. "0.999", "41.03") are stored as text, just
# like the original inlineStr cells, instead of being auto-converted to
# numbers. The Style is reset to "Normal" right after so that no stray
# "quote prefix" cell formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '''67.100.99'
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = '''  +2.53%  '
$ws.Cells.Item(2, 5).Style = "Normal"
# Row 3
$ws.Cells.Item(3, 4).Value = '''3.828.02'
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = '''  +4.04%  '
$ws.Cells.Item(3, 5).Style = "Normal"
# Row 4
$ws.Cells.Item(4, 4).Value = '''0.999'
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '''  -0.10%  '
$ws.Cells.Item(4, 5).Style = "Normal"
# Row 5
$ws.Cells.Item(5, 4).Value = '''423.16'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '''  +3.84%  '
$ws.Cells.Item(5, 5).Style = "Normal"
# Row 6
$ws.Cells.Item(6, 5).Value = '''  -3.71%  '
$ws.Cells.Item(6, 5).Style = "Normal"
# Row 7
$ws.Cells.Item(7, 4).Value = '''3.822.56'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '''  +4.07%  '
$ws.Cells.Item(7, 5).Style = "Normal"
# Row 8
$ws.Cells.Item(8, 4).Value = '''0.609'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '''  -2.45%  '
$ws.Cells.Item(8, 5).Style = "Normal"
# Row 9
$ws.Cells.Item(9, 5).Value = '''  -0.14%  '
$ws.Cells.Item(9, 5).Style = "Normal"
# Row 10
$ws.Cells.Item(10, 5).Value = '''  -1.48%  '
$ws.Cells.Item(10, 5).Style = "Normal"
# Row 11
$ws.Cells.Item(11, 5).Value = '''  -2.29%  '
$ws.Cells.Item(11, 5).Style = "Normal"
# Row 12
$ws.Cells.Item(12, 4).Value = '''0.0000348'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '''  +6.12%  '
$ws.Cells.Item(12, 5).Style = "Normal"
# Row 13
$ws.Cells.Item(13, 4).Value = '''40.89'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '''  -3.85%  '
$ws.Cells.Item(13, 5).Style = "Normal"
# Row 14
$ws.Cells.Item(14, 4).Value = '''10.17'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '''  +1.22%  '
$ws.Cells.Item(14, 5).Style = "Normal"
# Row 15
$ws.Cells.Item(15, 4).Value = '''4.423.27'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '''  +4.13%  '
$ws.Cells.Item(15, 5).Style = "Normal"
# Row 16
$ws.Cells.Item(16, 4).Value = '''15.71'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '''  +16.12%  '
$ws.Cells.Item(16, 5).Style = "Normal"
# Row 17
$ws.Cells.Item(17, 5).Value = '''  -0.59%  '
$ws.Cells.Item(17, 5).Style = "Normal"
# Row 18
$ws.Cells.Item(18, 4).Value = '''3.821.56'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '''  +4.03%  '
$ws.Cells.Item(18, 5).Style = "Normal"
# Row 19
$ws.Cells.Item(19, 4).Value = '''19.69'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '''  -2.03%  '
$ws.Cells.Item(19, 5).Style = "Normal"
# Row 20
$ws.Cells.Item(20, 4).Value = '''67.205.89'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '''  +2.76%  '
$ws.Cells.Item(20, 5).Style = "Normal"
# Row 21
$ws.Cells.Item(21, 5).Value = '''  -0.72%  '
$ws.Cells.Item(21, 5).Style = "Normal"
# Row 22
$ws.Cells.Item(22, 4).Value = '''408.59'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '''  -3.80%  '
$ws.Cells.Item(22, 5).Style = "Normal"
# Row 23
$ws.Cells.Item(23, 4).Value = '''15.06'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '''  -1.75%  '
$ws.Cells.Item(23, 5).Style = "Normal"
# Row 24
$ws.Cells.Item(24, 4).Value = '''84.19'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '''  -2.55%  '
$ws.Cells.Item(24, 5).Style = "Normal"
# Row 25
$ws.Cells.Item(25, 4).Value = '''3.05'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '''  +1.18%  '
$ws.Cells.Item(25, 5).Style = "Normal"
# Row 26
$ws.Cells.Item(26, 4).Value = '''37.15'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '''  +2.91%  '
$ws.Cells.Item(26, 5).Style = "Normal"
# Row 27
$ws.Cells.Item(27, 4).Value = '''10.06'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '''  +5.99%  '
$ws.Cells.Item(27, 5).Style = "Normal"
# Row 28
$ws.Cells.Item(28, 5).Value = '''  +0.88%  '
$ws.Cells.Item(28, 5).Style = "Normal"
# Row 29
$ws.Cells.Item(29, 4).Value = '''5.45'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '''  +6.14%  '
$ws.Cells.Item(29, 5).Style = "Normal"
# Row 30
$ws.Cells.Item(30, 4).Value = '''9.04'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '''  +28.85%  '
$ws.Cells.Item(30, 5).Style = "Normal"
# Row 31
$ws.Cells.Item(31, 4).Value = '''733.89'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '''  +6.53%  '
$ws.Cells.Item(31, 5).Style = "Normal"
# Row 32
$ws.Cells.Item(32, 4).Value = '''12.73'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '''  -0.59%  '
$ws.Cells.Item(32, 5).Style = "Normal"
# Row 33
$ws.Cells.Item(33, 4).Value = '''2.77'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '''  +1.89%  '
$ws.Cells.Item(33, 5).Style = "Normal"
# Row 34
$ws.Cells.Item(34, 5).Value = '''  +2.47%  '
$ws.Cells.Item(34, 5).Style = "Normal"
# Row 35
$ws.Cells.Item(35, 5).Value = '''  -0.05%  '
$ws.Cells.Item(35, 5).Style = "Normal"
# Row 36
$ws.Cells.Item(36, 5).Value = '''  -6.75%  '
$ws.Cells.Item(36, 5).Style = "Normal"
# Row 37
$ws.Cells.Item(37, 4).Value = '''38.69'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '''  -6.26%  '
$ws.Cells.Item(37, 5).Style = "Normal"
# Row 38
$ws.Cells.Item(38, 4).Value = '''55.37'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '''  -1.02%  '
$ws.Cells.Item(38, 5).Style = "Normal"
# Row 39
$ws.Cells.Item(39, 4).Value = '''5.46'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '''  +24.29%  '
$ws.Cells.Item(39, 5).Style = "Normal"
# Row 40
$ws.Cells.Item(40, 4).Value = '''0.0₃0753'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '''  +15.31%  '
$ws.Cells.Item(40, 5).Style = "Normal"
# Row 41
$ws.Cells.Item(41, 4).Value = '''0.0456'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '''  -2.32%  '
$ws.Cells.Item(41, 5).Style = "Normal"
# Row 42
$ws.Cells.Item(42, 4).Value = '''2.91'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '''  -0.32%  '
$ws.Cells.Item(42, 5).Style = "Normal"
# Row 43
$ws.Cells.Item(43, 5).Value = '''  +0.48%  '
$ws.Cells.Item(43, 5).Style = "Normal"
# Row 44
$ws.Cells.Item(44, 2).Value = '''LidoDAOToken'
$ws.Cells.Item(44, 2).Style = "Normal"
$ws.Cells.Item(44, 3).Value = '''https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(44, 3).Style = "Normal"
$ws.Cells.Item(44, 4).Value = '''3.35'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '''  +0.17%  '
$ws.Cells.Item(44, 5).Style = "Normal"
# Row 45
$ws.Cells.Item(45, 2).Value = '''Stellar'
$ws.Cells.Item(45, 2).Style = "Normal"
$ws.Cells.Item(45, 3).Value = '''https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(45, 3).Style = "Normal"
$ws.Cells.Item(45, 4).Value = '''0.134'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '''  -4.00%  '
$ws.Cells.Item(45, 5).Style = "Normal"
# Row 46
$ws.Cells.Item(46, 4).Value = '''3.15'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '''  +0.33%  '
$ws.Cells.Item(46, 5).Style = "Normal"
# Row 47
$ws.Cells.Item(47, 2).Value = '''Monero'
$ws.Cells.Item(47, 2).Style = "Normal"
$ws.Cells.Item(47, 3).Value = '''https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(47, 3).Style = "Normal"
$ws.Cells.Item(47, 4).Value = '''143.28'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '''  -0.70%  '
$ws.Cells.Item(47, 5).Style = "Normal"
# Row 48
$ws.Cells.Item(48, 2).Value = '''TheGraph'
$ws.Cells.Item(48, 2).Style = "Normal"
$ws.Cells.Item(48, 3).Value = '''https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(48, 3).Style = "Normal"
$ws.Cells.Item(48, 4).Value = '''0.313'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '''  +7.13%  '
$ws.Cells.Item(48, 5).Style = "Normal"
# Row 49
$ws.Cells.Item(49, 5).Value = '''  -2.05%  '
$ws.Cells.Item(49, 5).Style = "Normal"
# Row 50
$ws.Cells.Item(50, 4).Value = '''2.83'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '''  +0.43%  '
$ws.Cells.Item(50, 5).Style = "Normal"
# Row 51
$ws.Cells.Item(51, 4).Value = '''25.56'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '''  -5.45%  '
$ws.Cells.Item(51, 5).Style = "Normal"
